# DRILL-8417: Allow Excel Reader to Ignore Formula Errors
# Adds a new worksheet "Sheet with Errors" (after Sheet1) containing a
# field1/field2/result table whose "result" column divides field1 by
# field2 -- including a row that divides by zero and produces #DIV/0!.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Sheet1's selection moves to C4 (it is no longer the active/selected sheet
# once the new sheet is added).
$ws1.Range("C4").Select()

# Insert the new sheet right after Sheet1.
$ws = $wb.Worksheets.Add($null, $ws1)
$ws.Name = "Sheet with Errors"

# Header row -- set B1 before A1 so the shared-string table gets the
# same insertion order as the target workbook (field2, then field1).
$ws.Range("B1").Value = "field2"
$ws.Range("A1").Value = "field1"
$ws.Range("C1").Value = "result"

# Row 2: plain (non-shared) formula.
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 2
$ws.Range("C2").Formula = "=A2/B2"

# Rows 3-6: data first, then one shared formula written across the whole
# range so it collapses into a single shared-formula group.
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 3
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 4
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 0
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 6

$ws.Range("C3:C6").Formula = "=A3/B3"

# Final selection on the new (now active) sheet.
$ws.Range("E5").Select()
